$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) stays text, preserving formats like "1.020" or "27.214.51"
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '27.214.51'
$ws.Cells.Item(2, 5).Value = '  +1.14%  '
$ws.Cells.Item(3, 4).Value = '1.858.79'
$ws.Cells.Item(3, 5).Value = '  +0.72%  '
$ws.Cells.Item(4, 5).Value = '  +1.56%  '
$ws.Cells.Item(5, 5).Value = '  +1.27%  '
$ws.Cells.Item(6, 4).Value = '311.69'
$ws.Cells.Item(6, 5).Value = '  +0.68%  '
$ws.Cells.Item(7, 4).Value = '0.4791'
$ws.Cells.Item(7, 5).Value = '  +1.94%  '
$ws.Cells.Item(8, 4).Value = '0.3720'
$ws.Cells.Item(8, 5).Value = '  +1.51%  '
$ws.Cells.Item(9, 4).Value = '0.07308'
$ws.Cells.Item(9, 5).Value = '  +2.04%  '
$ws.Cells.Item(10, 4).Value = '0.9363'
$ws.Cells.Item(10, 5).Value = '  +0.88%  '
$ws.Cells.Item(11, 4).Value = '20.21'
$ws.Cells.Item(11, 5).Value = '  +3.19%  '
$ws.Cells.Item(12, 4).Value = '0.07871'
$ws.Cells.Item(12, 5).Value = '  +2.23%  '
$ws.Cells.Item(13, 4).Value = '1.848.33'
$ws.Cells.Item(13, 5).Value = '  -1.68%  '
$ws.Cells.Item(14, 4).Value = '5.420'
$ws.Cells.Item(14, 5).Value = '  +2.49%  '
$ws.Cells.Item(15, 4).Value = '6.540'
$ws.Cells.Item(15, 5).Value = '  +2.16%  '
$ws.Cells.Item(16, 4).Value = '90.30'
$ws.Cells.Item(16, 5).Value = '  +2.17%  '
$ws.Cells.Item(17, 5).Value = '  +1.26%  '
$ws.Cells.Item(18, 4).Value = '0.000008747'
$ws.Cells.Item(18, 5).Value = '  +1.22%  '
$ws.Cells.Item(19, 5).Value = '  +1.29%  '
$ws.Cells.Item(20, 4).Value = '27.251.61'
$ws.Cells.Item(20, 5).Value = '  +1.18%  '
$ws.Cells.Item(21, 4).Value = '14.71'
$ws.Cells.Item(22, 4).Value = '5.096'
$ws.Cells.Item(22, 5).Value = '  +1.50%  '
$ws.Cells.Item(23, 4).Value = '10.65'
$ws.Cells.Item(23, 5).Value = '  +0.37%  '
$ws.Cells.Item(24, 4).Value = '1.956'
$ws.Cells.Item(24, 5).Value = '  +1.29%  '
$ws.Cells.Item(25, 4).Value = '153.67'
$ws.Cells.Item(25, 5).Value = '  +1.15%  '
$ws.Cells.Item(26, 4).Value = '18.51'
$ws.Cells.Item(26, 5).Value = '  +1.43%  '
$ws.Cells.Item(27, 4).Value = '1.994'
$ws.Cells.Item(27, 5).Value = '  -0.96%  '
$ws.Cells.Item(28, 4).Value = '115.68'
$ws.Cells.Item(28, 5).Value = '  +1.11%  '
$ws.Cells.Item(29, 4).Value = '4.939'
$ws.Cells.Item(29, 5).Value = '  +1.21%  '
$ws.Cells.Item(30, 4).Value = '0.08885'
$ws.Cells.Item(30, 5).Value = '  +0.45%  '
$ws.Cells.Item(31, 4).Value = '3.345'
$ws.Cells.Item(31, 5).Value = '  +4.09%  '
$ws.Cells.Item(32, 5).Value = '  +0.14%  '
$ws.Cells.Item(33, 4).Value = '4.588'
$ws.Cells.Item(33, 5).Value = '  +2.49%  '
$ws.Cells.Item(34, 4).Value = '0.7405'
$ws.Cells.Item(34, 5).Value = '  -1.14%  '
$ws.Cells.Item(35, 4).Value = '2.678'
$ws.Cells.Item(35, 5).Value = '  -3.60%  '
$ws.Cells.Item(36, 5).Value = '  +3.58%  '
$ws.Cells.Item(37, 4).Value = '0.02020'
$ws.Cells.Item(37, 5).Value = '  +4.01%  '
$ws.Cells.Item(38, 2).Value = 'MXToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(38, 4).Value = '2.998'
$ws.Cells.Item(38, 5).Value = '  +1.49%  '
$ws.Cells.Item(39, 2).Value = 'Hedera'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(39, 4).Value = '0.05254'
$ws.Cells.Item(39, 5).Value = '  +0.79%  '
$ws.Cells.Item(40, 2).Value = 'TheSandbox'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(40, 4).Value = '0.5339'
$ws.Cells.Item(40, 5).Value = '  +2.06%  '
$ws.Cells.Item(41, 2).Value = 'FraxShare'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(41, 4).Value = '7.104'
$ws.Cells.Item(41, 5).Value = '  +1.78%  '
$ws.Cells.Item(42, 2).Value = 'Algorand'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(42, 4).Value = '0.1526'
$ws.Cells.Item(42, 5).Value = '  +0.88%  '
$ws.Cells.Item(43, 2).Value = 'Aptos'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(43, 4).Value = '8.335'
$ws.Cells.Item(43, 5).Value = '  +2.10%  '
$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(44, 4).Value = '10.57'
$ws.Cells.Item(44, 5).Value = '  +0.28%  '
$ws.Cells.Item(45, 2).Value = 'Decentraland'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(45, 4).Value = '0.4781'
$ws.Cells.Item(45, 5).Value = '  +1.50%  '
$ws.Cells.Item(46, 2).Value = 'PaxDollar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(46, 4).Value = '1.019'
$ws.Cells.Item(46, 5).Value = '  +1.23%  '
$ws.Cells.Item(47, 2).Value = 'Quant'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(47, 4).Value = '102.60'
$ws.Cells.Item(47, 5).Value = '  +1.71%  '
$ws.Cells.Item(48, 2).Value = 'NEARProtocol'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(48, 4).Value = '1.634'
$ws.Cells.Item(48, 5).Value = '  +2.20%  '
$ws.Cells.Item(49, 2).Value = 'Aave'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(49, 4).Value = '66.42'
$ws.Cells.Item(49, 5).Value = '  +1.66%  '
$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(50, 4).Value = '0.06076'
$ws.Cells.Item(50, 5).Value = '  +0.60%  '
$ws.Cells.Item(51, 2).Value = 'EOS'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Cells.Item(51, 4).Value = '0.8977'
$ws.Cells.Item(51, 5).Value = '  +0.49%  '
